$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1911.174
$ws.Range("I112").Value = 1024.75
$ws.Range("J112").Value = 2097.7896
$ws.Range("K112").Value = 3074.25
$ws.Range("L112").Value = 6293.3688
$ws.Range("M112").Value = -1966.25
$ws.Range("N112").Value = -8509.3688
$ws.Range("H125").Value = 995
$ws.Range("I125").Value = 990
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 8910
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -6450
$ws.Range("N125").Value = -13920
$ws.Range("H135").Value = 933.9259
$ws.Range("I135").Value = 655.2727
$ws.Range("J135").Value = 2160
$ws.Range("K135").Value = 5897.454299999999
$ws.Range("L135").Value = 19440
$ws.Range("M135").Value = -3362.454299999999
$ws.Range("N135").Value = -24510
$ws.Range("H138").Value = 3705.8086
$ws.Range("I138").Value = 1721
$ws.Range("J138").Value = 4112.9487
$ws.Range("K138").Value = 5163
$ws.Range("L138").Value = 12338.8461
$ws.Range("M138").Value = -23
$ws.Range("N138").Value = -22618.8461
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 981.2
$ws.Range("I45").Value = 801.3333
$ws.Range("J45").Value = 1147.2307
$ws.Range("K45").Value = 801.3333
$ws.Range("L45").Value = 1147.2307
$ws.Range("M45").Value = -424.3333
$ws.Range("N45").Value = -1901.2307
$ws.Range("H132").Value = 2939.4167
$ws.Range("I132").Value = 1530.8096
$ws.Range("J132").Value = 4911.467
$ws.Range("K132").Value = 4592.4288
$ws.Range("L132").Value = 14734.401
$ws.Range("M132").Value = -2062.4288
$ws.Range("N132").Value = -19794.401
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H58").Value = 2266.1865
$ws.Range("I58").Value = 1807.2
$ws.Range("J58").Value = 8577.25
$ws.Range("K58").Value = 1807.2
$ws.Range("L58").Value = 8577.25
$ws.Range("M58").Value = -1604.2
$ws.Range("N58").Value = -8983.25
$ws.Range("H59").Value = 25021
$ws.Range("J59").Value = 25021
$ws.Range("L59").Value = 25021
$ws.Range("N59").Value = -27311
$ws.Range("H99").Value = 4114.615
$ws.Range("I99").Value = 1866.6666
$ws.Range("J99").Value = 9172.5
$ws.Range("K99").Value = 1866.6666
$ws.Range("L99").Value = 9172.5
$ws.Range("M99").Value = -368.6666
$ws.Range("N99").Value = -12168.5
$ws.Range("H105").Value = 2130.35
$ws.Range("I105").Value = 800.5833
$ws.Range("K105").Value = 800.5833
$ws.Range("M105").Value = 946.4167
$ws.Range("H107").Value = 748.7857
$ws.Range("I107").Value = 667.9231
$ws.Range("K107").Value = 667.9231
$ws.Range("M107").Value = 1252.0769
$ws.Range("H126").Value = 4114.615
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 9172.5
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 27517.5
$ws.Range("M126").Value = -3129.9998
$ws.Range("N126").Value = -32457.5
$ws.Range("H134").Value = 4639.65
$ws.Range("I134").Value = 4507.3716
$ws.Range("K134").Value = 13522.1148
$ws.Range("M134").Value = -10987.1148
$ws.Range("H136").Value = 2266.1865
$ws.Range("I136").Value = 1807.2
$ws.Range("J136").Value = 8577.25
$ws.Range("K136").Value = 5421.6
$ws.Range("L136").Value = 25731.75
$ws.Range("M136").Value = -2871.6
$ws.Range("N136").Value = -30831.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1004964.8
$ws.Range("I4").Value = 1507034.8
$ws.Range("J4").Value = 825
$ws.Range("K4").Value = 4521104.4
$ws.Range("L4").Value = 2475
$ws.Range("M4").Value = -4520992.4
$ws.Range("N4").Value = -2699
$ws.Range("H107").Value = 1211.9
$ws.Range("I107").Value = 475
$ws.Range("J107").Value = 1948.8
$ws.Range("K107").Value = 1425
$ws.Range("L107").Value = 5846.4
$ws.Range("M107").Value = 495
$ws.Range("N107").Value = -9686.4
$ws.Range("H113").Value = 587.6111
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 587.6111
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1762.8333
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6102.8333
$ws.Range("H133").Value = 3896
$ws.Range("I133").Value = 5140
$ws.Range("J133").Value = 2807.5
$ws.Range("K133").Value = 15420
$ws.Range("L133").Value = 8422.5
$ws.Range("M133").Value = -10360
$ws.Range("N133").Value = -18542.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5693.4443
$ws.Range("I7").Value = 3082.3333
$ws.Range("J7").Value = 10915.667
$ws.Range("K7").Value = 3082.3333
$ws.Range("L7").Value = 10915.667
$ws.Range("M7").Value = -2970.3333
$ws.Range("N7").Value = -11139.667
$ws.Range("H13").Value = 14477.8
$ws.Range("I13").Value = 13990
$ws.Range("J13").Value = 14599.75
$ws.Range("K13").Value = 13990
$ws.Range("L13").Value = 14599.75
$ws.Range("M13").Value = -13850
$ws.Range("N13").Value = -14879.75
$ws.Range("H40").Value = 5491.647
$ws.Range("I40").Value = 4836.5713
$ws.Range("J40").Value = 8548.666999999999
$ws.Range("K40").Value = 4836.5713
$ws.Range("L40").Value = 8548.666999999999
$ws.Range("M40").Value = -4700.5713
$ws.Range("N40").Value = -8820.666999999999
$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H88").Value = 4723.6665
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H91").Value = 4723.6665
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H126").Value = 5693.4443
$ws.Range("I126").Value = 3082.3333
$ws.Range("J126").Value = 10915.667
$ws.Range("K126").Value = 9246.999899999999
$ws.Range("L126").Value = 32747.001
$ws.Range("M126").Value = -6776.999899999999
$ws.Range("N126").Value = -37687.001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 55680690
$ws.Range("I62").Value = 125003590
$ws.Range("J62").Value = 222361.2
$ws.Range("K62").Value = 125003590
$ws.Range("L62").Value = 222361.2
$ws.Range("M62").Value = -125002966
$ws.Range("N62").Value = -223609.2
$ws.Range("H65").Value = 55680690
$ws.Range("I65").Value = 125003590
$ws.Range("J65").Value = 222361.2
$ws.Range("K65").Value = 625017950
$ws.Range("L65").Value = 1111806
$ws.Range("M65").Value = -625014830
$ws.Range("N65").Value = -1118046
$ws.Range("H92").Value = 39800
$ws.Range("J92").Value = 39800
$ws.Range("L92").Value = 39800
$ws.Range("N92").Value = -44792
$ws.Range("H132").Value = 8551560
$ws.Range("I132").Value = 5634.381
$ws.Range("J132").Value = 18521806
$ws.Range("K132").Value = 16903.143
$ws.Range("L132").Value = 55565418
$ws.Range("M132").Value = -14373.143
$ws.Range("N132").Value = -55570478
